$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2996.818
$ws.Range("I100").Value = 2905
$ws.Range("J100").Value = 3006
$ws.Range("K100").Value = 2905
$ws.Range("L100").Value = 3006
$ws.Range("M100").Value = -2364
$ws.Range("N100").Value = -4088

$ws.Range("H118").Value = 451.26666
$ws.Range("I118").Value = 235
$ws.Range("J118").Value = 698.4286
$ws.Range("K118").Value = 705
$ws.Range("L118").Value = 2095.2858
$ws.Range("M118").Value = 952
$ws.Range("N118").Value = -5409.2858

$ws.Range("H125").Value = 9582.182000000001
$ws.Range("I125").Value = 258
$ws.Range("J125").Value = 20771.2
$ws.Range("K125").Value = 2322
$ws.Range("L125").Value = 186940.8
$ws.Range("M125").Value = 138
$ws.Range("N125").Value = -191860.8

$ws.Range("H138").Value = 1688901.5
$ws.Range("I138").Value = 4430.1113
$ws.Range("J138").Value = 2025795.8
$ws.Range("K138").Value = 13290.3339
$ws.Range("L138").Value = 6077387.4
$ws.Range("M138").Value = -8150.333899999998
$ws.Range("N138").Value = -6087667.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 930.5
$ws.Range("I4").Value = 611
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 611
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -495
$ws.Range("N4").Value = -1482

$ws.Range("H9").Value = 26002.25
$ws.Range("J9").Value = 26002.25
$ws.Range("L9").Value = 26002.25
$ws.Range("N9").Value = -26342.25

$ws.Range("H20").Value = 26002.25
$ws.Range("J20").Value = 26002.25
$ws.Range("L20").Value = 26002.25
$ws.Range("N20").Value = -26542.25

$ws.Range("H23").Value = 42003.6
$ws.Range("J23").Value = 36668.668
$ws.Range("L23").Value = 36668.668
$ws.Range("N23").Value = -37186.668

$ws.Range("H45").Value = 1511.1875
$ws.Range("I45").Value = 1525.24
$ws.Range("J45").Value = 1461
$ws.Range("K45").Value = 1525.24
$ws.Range("L45").Value = 1461
$ws.Range("M45").Value = -1148.24
$ws.Range("N45").Value = -2215

$ws.Range("H61").Value = 5882.2
$ws.Range("I61").Value = 2873.0698
$ws.Range("K61").Value = 2873.0698
$ws.Range("M61").Value = -2661.0698

$ws.Range("H74").Value = 3091.3115
$ws.Range("I74").Value = 1410.4889
$ws.Range("J74").Value = 7818.625
$ws.Range("K74").Value = 1410.4889
$ws.Range("L74").Value = 7818.625
$ws.Range("M74").Value = -536.4889000000001
$ws.Range("N74").Value = -9566.625

$ws.Range("H77").Value = 3091.3115
$ws.Range("I77").Value = 1410.4889
$ws.Range("J77").Value = 7818.625
$ws.Range("K77").Value = 7052.444500000001
$ws.Range("L77").Value = 39093.125
$ws.Range("M77").Value = -2684.444500000001
$ws.Range("N77").Value = -47829.125

$ws.Range("H88").Value = 5083.222
$ws.Range("I88").Value = 10220.667
$ws.Range("J88").Value = 2514.5
$ws.Range("K88").Value = 10220.667
$ws.Range("L88").Value = 2514.5
$ws.Range("M88").Value = -9814.666999999999
$ws.Range("N88").Value = -3326.5

$ws.Range("H91").Value = 5083.222
$ws.Range("I91").Value = 10220.667
$ws.Range("J91").Value = 2514.5
$ws.Range("K91").Value = 10220.667
$ws.Range("L91").Value = 2514.5
$ws.Range("M91").Value = -8816.666999999999
$ws.Range("N91").Value = -5322.5

$ws.Range("H136").Value = 5882.2
$ws.Range("I136").Value = 2873.0698
$ws.Range("K136").Value = 8619.2094
$ws.Range("M136").Value = -6069.2094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2767.65
$ws.Range("I86").Value = 2868.5293
$ws.Range("J86").Value = 2196
$ws.Range("K86").Value = 2868.5293
$ws.Range("L86").Value = 2196
$ws.Range("M86").Value = -1745.5293
$ws.Range("N86").Value = -4442

$ws.Range("H89").Value = 2767.65
$ws.Range("I89").Value = 2868.5293
$ws.Range("J89").Value = 2196
$ws.Range("K89").Value = 14342.6465
$ws.Range("L89").Value = 10980
$ws.Range("M89").Value = -8726.646500000001
$ws.Range("N89").Value = -22212

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 2815
$ws.Range("I94").Value = 2815
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2815
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2364
$ws.Range("N94").ClearContents()

$ws.Range("H104").Value = 59999.5
$ws.Range("J104").Value = 59999.5
$ws.Range("L104").Value = 59999.5
$ws.Range("N104").Value = -66987.5

$ws.Range("H107").Value = 1356.619
$ws.Range("I107").Value = 1321.625
$ws.Range("K107").Value = 1321.625
$ws.Range("M107").Value = 598.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 2658487.8
$ws.Range("I57").Value = 12000
$ws.Range("J57").Value = 2923136.5
$ws.Range("K57").Value = 12000
$ws.Range("L57").Value = 2923136.5
$ws.Range("M57").Value = -11440
$ws.Range("N57").Value = -2924256.5

$ws.Range("H58").Value = 1654868.9
$ws.Range("I58").Value = 2067574
$ws.Range("J58").Value = 4048.5454
$ws.Range("K58").Value = 2067574
$ws.Range("L58").Value = 4048.5454
$ws.Range("M58").Value = -2067371
$ws.Range("N58").Value = -4454.5454

$ws.Range("H122").Value = 14694.389
$ws.Range("I122").Value = 7826.5
$ws.Range("J122").Value = 20188.7
$ws.Range("K122").Value = 23479.5
$ws.Range("L122").Value = 60566.10000000001
$ws.Range("M122").Value = -21029.5
$ws.Range("N122").Value = -65466.10000000001

$ws.Range("H136").Value = 1654868.9
$ws.Range("I136").Value = 2067574
$ws.Range("J136").Value = 4048.5454
$ws.Range("K136").Value = 6202722
$ws.Range("L136").Value = 12145.6362
$ws.Range("M136").Value = -6200172
$ws.Range("N136").Value = -17245.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2780771.5
$ws.Range("J5").Value = 5560860.5
$ws.Range("L5").Value = 16682581.5
$ws.Range("N5").Value = -16682805.5

$ws.Range("H69").Value = 14287025
$ws.Range("J69").Value = 15626387
$ws.Range("L69").Value = 46879161
$ws.Range("N69").Value = -46880783

$ws.Range("H72").Value = 14287025
$ws.Range("J72").Value = 15626387
$ws.Range("L72").Value = 140637483
$ws.Range("N72").Value = -140645595

$ws.Range("H122").Value = 1470.4706
$ws.Range("I122").Value = 664.4
$ws.Range("J122").Value = 1806.3334
$ws.Range("K122").Value = 5979.599999999999
$ws.Range("L122").Value = 16257.0006
$ws.Range("M122").Value = -3529.599999999999
$ws.Range("N122").Value = -21157.0006

$ws.Range("H129").Value = 2306.25
$ws.Range("I129").Value = 2945.5557
$ws.Range("J129").Value = 1484.2858
$ws.Range("K129").Value = 8836.667099999999
$ws.Range("L129").Value = 4452.857400000001
$ws.Range("M129").Value = -3836.667099999999
$ws.Range("N129").Value = -14452.8574

$ws.Range("H131").Value = 10260.55
$ws.Range("I131").Value = 476.3889
$ws.Range("J131").Value = 21746.305
$ws.Range("K131").Value = 1429.1667
$ws.Range("L131").Value = 65238.915
$ws.Range("M131").Value = 3610.8333
$ws.Range("N131").Value = -75318.91500000001

$ws.Range("H132").Value = 1581.1428
$ws.Range("I132").Value = 1329.4117
$ws.Range("J132").Value = 2651
$ws.Range("K132").Value = 11964.7053
$ws.Range("L132").Value = 23859
$ws.Range("M132").Value = -9434.705300000001
$ws.Range("N132").Value = -28919

$ws.Range("H135").Value = 2780771.5
$ws.Range("J135").Value = 5560860.5
$ws.Range("L135").Value = 50047744.5
$ws.Range("N135").Value = -50052814.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 263.9091
$ws.Range("I2").Value = 101.25
$ws.Range("J2").Value = 356.85715
$ws.Range("K2").Value = 101.25
$ws.Range("L2").Value = 356.85715
$ws.Range("M2").Value = 11.75
$ws.Range("N2").Value = -582.85715

$ws.Range("H102").Value = 6340.8237
$ws.Range("I102").Value = 5890.909
$ws.Range("J102").Value = 7165.6665
$ws.Range("K102").Value = 5890.909
$ws.Range("L102").Value = 7165.6665
$ws.Range("M102").Value = -4268.909
$ws.Range("N102").Value = -10409.6665

$ws.Range("H140").Value = 59556
$ws.Range("J140").Value = 59556
$ws.Range("L140").Value = 59556
$ws.Range("N140").Value = -69916

$ws.Range("H141").Value = 68423
$ws.Range("J141").Value = 68423
$ws.Range("L141").Value = 68423
$ws.Range("N141").Value = -78783

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9666.666999999999
$ws.Range("I93").Value = 10000
$ws.Range("J93").Value = 9000
$ws.Range("K93").Value = 10000
$ws.Range("L93").Value = 9000
$ws.Range("M93").Value = -8752
$ws.Range("N93").Value = -11496

$ws.Range("H132").Value = 6883.4707
$ws.Range("I132").Value = 9301.200000000001
$ws.Range("J132").Value = 3429.5715
$ws.Range("K132").Value = 27903.6
$ws.Range("L132").Value = 10288.7145
$ws.Range("M132").Value = -25373.6
$ws.Range("N132").Value = -15348.7145

$ws.Range("H136").Value = 6207.6
$ws.Range("I136").Value = 4359
$ws.Range("J136").Value = 8625
$ws.Range("K136").Value = 13077
$ws.Range("L136").Value = 25875
$ws.Range("M136").Value = -10527
$ws.Range("N136").Value = -30975

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4587

$ws.Range("H81").Value = 4959.3335
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4959.3335
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 9918.666999999999
$ws.Range("N81").Value = -12040.667
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 4959.3335
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4959.3335
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 49593.335
$ws.Range("N84").Value = -60201.335
$ws.Range("M84").ClearContents()

$ws.Range("H140").Value = 55813.332
$ws.Range("J140").Value = 55813.332
$ws.Range("L140").Value = 55813.332
$ws.Range("N140").Value = -66173.33199999999

$ws.Range("H141").Value = 45614.168
$ws.Range("J141").Value = 46974.547
$ws.Range("L141").Value = 46974.547
$ws.Range("N141").Value = -57334.547
